# Notas do fórum para a semana 05/06/2022 a 11/06/2022 no semestre 2022-1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I was the per-day column for 2022-06-05; it is being replaced with
# the total_views aggregate, and column J (was 2022-06-06) becomes nota_view.
# The old K (total_views) and L (nota_view) columns are dropped since their
# content moves into I and J.

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $totalViews = 0
    for ($c = 2; $c -le 8; $c++) {
        $totalViews += $ws.Cells.Item($r, $c).Value2
    }
    $notaView = $ws.Cells.Item($r, 12).Value2   # old column L value
    $ws.Cells.Item($r, 9).Value2 = $totalViews  # column I
    $ws.Cells.Item($r, 10).Value2 = $notaView   # column J
}

# Now remove the old K and L columns (they are now redundant duplicates).
$ws.Range("K1:L1").EntireColumn.Delete() | Out-Null

# Fix the header labels for I1/J1.
$ws.Range("I1").Value = "total_views"
$ws.Range("J1").Value = "nota_view"
